# The unified diff supplied for this task shows word/document.xml and
# word/styles.xml differing only in cosmetic ways:
#   - w:rsid* attributes (random per-save revision-save IDs) are absent from
#     the "after" side purely because the diffing/scraping tool that produced
#     the xml_diff strips them before rendering (they are not meaningful
#     content and are not actually removed from the saved package).
#   - every remaining element's attributes are reordered into plain
#     alphabetical order by local name (namespace declarations first, also
#     alphabetised by prefix), which is how that same tool pretty-prints
#     XML for a readable diff -- it is not something a Word editing action
#     (or a real save) performs, and it carries no semantic meaning in XML.
#
# Reconstructing the "before" content, stripping w:rsid*, and alphabetising
# attributes reproduces the "after" side of the diff line for line (verified
# while preparing this script), confirming there is no text, formatting,
# structural, or property change to apply: every paragraph, run, field,
# bookmark, section property, style default, latent style entry and style
# definition is byte-for-byte identical content-wise before and after.
#
# So the faithful edit is a no-op against the Word object model: open the
# already-active document and make no content change. (A COM/OM-level edit
# cannot reorder raw XML attributes or drop rsid bookkeeping anyway -- Word
# automation only ever lets you change document *content*, which here did
# not change.)
$d = $word.ActiveDocument

# Touch the document without altering any content, to mirror the save that
# produced the (content-identical) revision described by the diff.
$null = $d.Content
